# sm_car_data_Tire_Tire2x.xlsx — "Update 2p0. Convention change to support
# multi-axle vehicles"
#
# The single-axle sheet "Bus_Makhulu_2x" is replaced by two tire sheets that
# follow the new naming convention:
#   - Tire2x_270_70R22  (keeps the original tab's data / position)
#   - Tire2x_430_50R38  (a new tab, duplicated from the first, becomes active)
#
# Both tabs keep the same layout/instance name ("Tire") in H2, but each gets
# its own Tire2x instance name in H3, and the xOffset value in H7 is stored
# as a plain (no longer formula-derived) value of 0.4572 m.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the existing sheet so both tabs share identical layout, column
# widths, styles and conditional formatting.
$ws1.Copy($null, $ws1)

$ws270 = $wb.Worksheets.Item(1)
$ws430 = $wb.Worksheets.Item(2)

$ws270.Name = "Tire2x_270_70R22"
$ws430.Name = "Tire2x_430_50R38"

# --- Tire2x_270_70R22 -------------------------------------------------
$ws270.Range("H2").Value = "Tire"
$ws270.Range("H3").Value = "Tire2x_270_70R22"
$ws270.Range("H7").Value = 0.4572
[void]$ws270.Range("C25").Select()

# --- Tire2x_430_50R38 --------------------------------------------------
$ws430.Range("H2").Value = "Tire"
$ws430.Range("H3").Value = "Tire2x_430_50R38"
$ws430.Range("H7").Value = 0.4572
[void]$ws430.Range("J16").Select()

# The newly added tab becomes the active/selected one.
$ws430.Activate()
